# Update gh-pages to output generated at 456a3b4
# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 293
$ws1.Range("F3").Value = 20
$ws1.Range("F4").Value = 7865
$ws1.Range("F5").Value = 5748
$ws1.Range("F6").Value = 481
$ws1.Range("F10").Value = 268
$ws1.Range("F11").Value = 303
$ws1.Range("F12").Value = 62

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 293
$ws4.Range("F3").Value = 20
$ws4.Range("F4").Value = 7865
$ws4.Range("F5").Value = 5748
$ws4.Range("F6").Value = 481
$ws4.Range("F10").Value = 268
$ws4.Range("F14").Value = 303
$ws4.Range("F15").Value = 62

$wb.Save()
